$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 14
$ws.Range("H2").Value = 14

$ws.Range("E7").Value = 29
$ws.Range("F7").Value = 16
$ws.Range("H7").Value = 16

$ws.Range("E8").Value = 42

$ws.Range("F14").Value = 18
$ws.Range("H14").Value = 18

$ws.Range("E15").Value = 92

$ws.Range("E16").Value = 304
